# Natmi following Dr Hou advice
# Update LR-pair rows for Serpinc1-Sdc2 (ligand-detection-rate threshold raised: 1 -> 3 cells),
# recomputing all derived statistics, and append new rows 8-10 for the sCs sending cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, "ECs", "Serpinc1", "Sdc2", "ECs", 3, 1, 4.117659, 12.352977, 0.411783455701328, 0.4117834557013281, 2, 0.6666666666666666, 1.321445333333333, 3.964336, 0.01021782062667047, 0.01021782062667047, 5.441261269808, 48.971351428272, 0.004207529487386674, 0.004207529487386676),
    @(3, "ECs", "Serpinc1", "Sdc2", "FAPs", 3, 1, 4.117659, 12.352977, 0.411783455701328, 0.4117834557013281, 3, 1, 105.9632263333333, 317.889679, 0.819340166699254, 0.8193401666992541, 436.320432580487, 3926.883893224383, 0.337390725238321, 0.3373907252383211),
    @(4, "ECs", "Serpinc1", "Sdc2", "sCs", 3, 1, 4.117659, 12.352977, 0.411783455701328, 0.4117834557013281, 3, 1, 22.04284166666666, 66.128525, 0.1704420126740755, 0.1704420126740755, 90.76490537432498, 816.8841483689249, 0.07018520097562034, 0.07018520097562037),
    @(5, "FAPs", "Serpinc1", "Sdc2", "ECs", 3, 1, 3.004148, 9.012444, 0.3004276082303642, 0.3004276082303643, 2, 0.6666666666666666, 1.321445333333333, 3.964336, 0.01021782062667047, 0.01021782062667047, 3.969817355242667, 35.72835619718401, 0.00306971541219749, 0.003069715412197491),
    @(6, "FAPs", "Serpinc1", "Sdc2", "FAPs", 3, 1, 3.004148, 9.012444, 0.3004276082303642, 0.3004276082303643, 3, 1, 105.9632263333333, 317.889679, 0.819340166699254, 0.8193401666992541, 318.3292144628307, 2864.962930165476, 0.2461524066085248, 0.2461524066085249),
    @(7, "FAPs", "Serpinc1", "Sdc2", "sCs", 3, 1, 3.004148, 9.012444, 0.3004276082303642, 0.3004276082303643, 3, 1, 22.04284166666666, 66.128525, 0.1704420126740755, 0.1704420126740755, 66.21995870723333, 595.9796283651, 0.05120548620964192, 0.05120548620964194),
    @(8, "sCs", "Serpinc1", "Sdc2", "ECs", 3, 1, 2.877766666666667, 8.6333, 0.2877889360683077, 0.2877889360683077, 2, 0.6666666666666666, 1.321445333333333, 3.964336, 0.01021782062667047, 0.01021782062667047, 3.802811332088889, 34.22530198880001, 0.002940575727086303, 0.002940575727086304),
    @(9, "sCs", "Serpinc1", "Sdc2", "FAPs", 3, 1, 2.877766666666667, 8.6333, 0.2877889360683077, 0.2877889360683077, 3, 1, 105.9632263333333, 317.889679, 0.819340166699254, 0.8193401666992541, 304.9374406345222, 2744.4369657107, 0.2357970348524082, 0.2357970348524082),
    @(10, "sCs", "Serpinc1", "Sdc2", "sCs", 3, 1, 2.877766666666667, 8.6333, 0.2877889360683077, 0.2877889360683077, 3, 1, 22.04284166666666, 66.128525, 0.1704420126740755, 0.1704420126740755, 63.43415498694443, 570.9073948825, 0.04905132548881319, 0.04905132548881321)
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($r, $col).Value2 = $row[$col]
    }
}
